# Updates the cryptocurrency price/volume table (columns D and E, rows 2-51)
# and fixes the swapped MXToken / WEMIXToken rows (B41:E42) per the refreshed feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is @(CellAddress, NewValue). A leading "'" forces a numeric-looking
# string (e.g. "38.00") to stay text instead of being coerced to a Number.
$edits = @(
    @('D2', '34.442.89'),
    @('E2', '  +0.28%  '),
    @('D3', '1.803.55'),
    @('E3', '  -0.06%  '),
    @('E4', '  +0.02%  '),
    @('D5', '''224.65'),
    @('E5', '  -1.36%  '),
    @('D6', '''0.588'),
    @('E6', '  +2.33%  '),
    @('E7', '  +0.01%  '),
    @('D8', '''38.00'),
    @('E8', '  +4.49%  '),
    @('E9', '  -5.10%  '),
    @('E10', '  -3.59%  '),
    @('D11', '''0.0973'),
    @('E11', '  +1.16%  '),
    @('D12', '2.064.81'),
    @('E12', '  -0.06%  '),
    @('D13', '''11.03'),
    @('E13', '  -6.96%  '),
    @('D14', '1.808.10'),
    @('E14', '  +0.34%  '),
    @('D15', '34.421.19'),
    @('E15', '  +0.20%  '),
    @('D16', '''0.624'),
    @('E16', '  -3.28%  '),
    @('E17', '  -3.45%  '),
    @('D18', '''67.74'),
    @('E18', '  -2.26%  '),
    @('D19', '''241.21'),
    @('E19', '  -2.02%  '),
    @('E20', '  -3.80%  '),
    @('D21', '''10.99'),
    @('E21', '  -6.55%  '),
    @('E22', '  +0.02%  '),
    @('D23', '''4.07'),
    @('E23', '  -2.84%  '),
    @('D24', '''2.17'),
    @('E24', '  +1.85%  '),
    @('D25', '''170.31'),
    @('E25', '  -0.77%  '),
    @('D26', '''7.68'),
    @('E26', '  -3.68%  '),
    @('D27', '''17.41'),
    @('E27', '  +2.90%  '),
    @('D28', '''0.119'),
    @('E28', '  +0.21%  '),
    @('E29', '  +0.03%  '),
    @('E30', '  -1.57%  '),
    @('D31', '''3.74'),
    @('E31', '  -2.84%  '),
    @('E32', '  -4.88%  '),
    @('E33', '  -3.93%  '),
    @('E34', '  -0.95%  '),
    @('D35', '1.318.91'),
    @('E35', '  -5.64%  '),
    @('E36', '  -5.54%  '),
    @('E37', '  -1.74%  '),
    @('E38', '  -2.14%  '),
    @('E39', '  +0.76%  '),
    @('D40', '''2.28'),
    @('E40', '  -8.33%  '),
    @('B41', 'WEMIXToken'),
    @('C41', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'),
    @('D41', '''1.21'),
    @('E41', '  -1.56%  '),
    @('B42', 'MXToken'),
    @('C42', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'),
    @('D42', '''2.80'),
    @('E42', '  -0.64%  '),
    @('D43', '''81.25'),
    @('E43', '  -1.76%  '),
    @('E44', '  -2.92%  '),
    @('D45', '''13.64'),
    @('E45', '  +1.59%  '),
    @('E46', '  +0.79%  '),
    @('D47', '1.965.36'),
    @('E47', '  -0.09%  '),
    @('E48', '  -5.59%  '),
    @('E49', '  +0.02%  '),
    @('D50', '''101.74'),
    @('E50', '  -2.76%  '),
    @('D51', '0.0₆0121'),
    @('E51', '  -5.86%  ')
)

foreach ($edit in $edits) {
    $ref = $edit[0]
    $value = $edit[1]
    $range = $ws.Range($ref)
    $range.Value = $value
    # Re-apply the default style so forcing text doesn't leave a stray number format behind.
    $range.Style = 'Normal'
}

Write-Host "Applied $($edits.Count) cell updates to $($ws.Name)"
